$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was: 05-11-2021, Californiana(o)/Primera, 100, bandeja) -> now 07-12-2022, Golden Nugget/Especial, 60, caja
$ws.Range("D2").Value = 44902
$ws.Range("K2").Value = "Golden Nugget"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "`$/caja 10 kilos"
$ws.Range("S2").Value = 1500

# Row 3 (was: 05-11-2021, Golden Nugget/Primera, 50, bandeja) -> now 07-12-2022, Golden Nugget/Primera, 70, caja
$ws.Range("D3").Value = 44902
$ws.Range("M3").Value = 70
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 13000
$ws.Range("Q3").Value = "`$/caja 10 kilos"
$ws.Range("S3").Value = 1300

# Row 4 (was: 07-12-2022, Golden Nugget/Especial, 60, caja) -> now 05-11-2021, Californiana(o)/Primera, 100, bandeja
$ws.Range("D4").Value = 44505
$ws.Range("K4").Value = "Californiana(o)"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "`$/bandeja 10 kilos"
$ws.Range("S4").Value = 1500

# Row 5 (was: 07-12-2022, Golden Nugget/Primera, 70, caja) -> now 05-11-2021, Golden Nugget/Primera, 50, bandeja
$ws.Range("D5").Value = 44505
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "`$/bandeja 10 kilos"
$ws.Range("S5").Value = 1500
